$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 615
$ws.Range("I55").Value = 615
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 615
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -401
$ws.Range("N55").ClearContents()

$ws.Range("H64").Value = 5000
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496

$ws.Range("H67").Value = 5000
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716

$ws.Range("H69").Value = 4250
$ws.Range("I69").Value = 2000
$ws.Range("K69").Value = 6000
$ws.Range("M69").Value = -5126

$ws.Range("H72").Value = 4250
$ws.Range("I72").Value = 2000
$ws.Range("K72").Value = 18000
$ws.Range("M72").Value = -13632

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H105").Value = 56225
$ws.Range("J105").Value = 56225
$ws.Range("L105").Value = 56225
$ws.Range("N105").Value = -63213

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 119999
$ws.Range("J107").Value = 119999
$ws.Range("L107").Value = 119999
$ws.Range("N107").Value = -127679

$ws.Range("H117").Value = 99999
$ws.Range("J117").Value = 99999
$ws.Range("L117").Value = 99999
$ws.Range("N117").Value = -109177

$ws.Range("H132").Value = 11903.333
$ws.Range("I132").Value = 14204.4
$ws.Range("J132").Value = 398
$ws.Range("K132").Value = 42613.2
$ws.Range("L132").Value = 1194
$ws.Range("M132").Value = -40083.2
$ws.Range("N132").Value = -6254

$ws.Range("H139").Value = 71000
$ws.Range("J139").Value = 71000
$ws.Range("L139").Value = 71000
$ws.Range("N139").Value = -81280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 264.66666
$ws.Range("I7").Value = 264.66666
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 264.66666
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -151.66666
$ws.Range("N7").ClearContents()

$ws.Range("H99").Value = 1200
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H103").Value = 17703.857
$ws.Range("J103").Value = 17703.857
$ws.Range("L103").Value = 17703.857
$ws.Range("N103").Value = -20047.857

$ws.Range("H106").Value = 14166
$ws.Range("J106").Value = 14166
$ws.Range("L106").Value = 14166
$ws.Range("N106").Value = -16690

$ws.Range("H112").Value = 119999
$ws.Range("J112").Value = 119999
$ws.Range("L112").Value = 119999
$ws.Range("N112").Value = -122953

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H4").Value = 5000
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5224

$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H28").Value = 33999
$ws.Range("J28").Value = 33999
$ws.Range("L28").Value = 33999
$ws.Range("N28").Value = -34489

$ws.Range("H31").Value = 1499.6666
$ws.Range("I31").Value = 1499.6666
$ws.Range("K31").Value = 1499.6666
$ws.Range("M31").Value = -1204.6666

$ws.Range("H34").Value = 1499.6666
$ws.Range("I34").Value = 1499.6666
$ws.Range("K34").Value = 1499.6666
$ws.Range("M34").Value = -1297.6666

$ws.Range("H43").Value = 11575.75
$ws.Range("J43").Value = 11575.75
$ws.Range("L43").Value = 11575.75
$ws.Range("N43").Value = -11943.75

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H74").Value = 15000
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 12500
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 12500
$ws.Range("M74").Value = -19126
$ws.Range("N74").Value = -14248

$ws.Range("H77").Value = 15000
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 12500
$ws.Range("K77").Value = 60000
$ws.Range("L77").Value = 37500
$ws.Range("M77").Value = -55632
$ws.Range("N77").Value = -46236

$ws.Range("H88").Value = 14779.5
$ws.Range("J88").Value = 14779.5
$ws.Range("L88").Value = 14779.5
$ws.Range("N88").Value = -15591.5

$ws.Range("H91").Value = 14779.5
$ws.Range("J91").Value = 14779.5
$ws.Range("L91").Value = 14779.5
$ws.Range("N91").Value = -17587.5

$ws.Range("H95").Value = 6500
$ws.Range("J95").Value = 6500
$ws.Range("L95").Value = 6500
$ws.Range("N95").Value = -11992

$ws.Range("H101").Value = 11575.75
$ws.Range("J101").Value = 11575.75
$ws.Range("L101").Value = 11575.75
$ws.Range("N101").Value = -18065.75

$ws.Range("H122").Value = 1157.25
$ws.Range("I122").Value = 1157.25
$ws.Range("K122").Value = 3471.75
$ws.Range("M122").Value = -1021.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 34999
$ws.Range("J98").Value = 34999
$ws.Range("L98").Value = 34999
$ws.Range("N98").Value = -40989

$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2251

$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256

$ws.Range("H82").Value = 1499.75

$ws.Range("H85").Value = 1499.75

$ws.Range("H93").Value = 1698.6666
$ws.Range("I93").Value = 1438.6
$ws.Range("K93").Value = 1438.6
$ws.Range("M93").Value = -190.5999999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3449.5
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3449.5
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 6899
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -7981

$ws.Range("H113").Value = 1266.3334
$ws.Range("I113").Value = 1266.3334
$ws.Range("K113").Value = 3799.0002
$ws.Range("M113").Value = -1629.0002

$ws.Range("H136").Value = 1768.4546
$ws.Range("I136").Value = 1768.4546
$ws.Range("K136").Value = 5305.3638
$ws.Range("M136").Value = -2755.3638

$ws.Range("H140").Value = 99883
$ws.Range("J140").Value = 99883
$ws.Range("L140").Value = 99883
$ws.Range("N140").Value = -110243

$ws.Range("H141").Value = 99495
$ws.Range("J141").Value = 99495
$ws.Range("L141").Value = 99495
$ws.Range("N141").Value = -109855
